# dataset update 17 april
# Add 2020-04-17 (serial 43938) data to all three sheets and update the
# active-sheet / selection bookkeeping that Excel rewrites on save.
#
# NB: Range.Select() implicitly activates the sheet it belongs to (same as
# real Excel), so the sheet whose selection should "win" as the final
# ActiveSheet / tabSelected / workbook.activeTab must be touched last.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Confirmed (sheet1) - new row 42
# ---------------------------------------------------------------------
$wsConfirmed = $wb.Worksheets.Item("Confirmed")

$wsConfirmed.Range("A42").Value = 43938
$wsConfirmed.Range("A42").NumberFormat = "d-mmm-yy"

$wsConfirmed.Range("C42").Value = 266
$wsConfirmed.Range("C42").HorizontalAlignment = -4108
$wsConfirmed.Range("C42").VerticalAlignment = -4108

$wsConfirmed.Range("B42").Formula = "=SUM(B41+C42)"
$wsConfirmed.Range("B42").HorizontalAlignment = -4108
$wsConfirmed.Range("B42").VerticalAlignment = -4108

# Confirmed's selection moves, and it loses tabSelected (done first, it
# will be superseded once another sheet is activated below).
$wsConfirmed.Range("B42").Select()

# ---------------------------------------------------------------------
# Death (sheet3) - new row 42
# ---------------------------------------------------------------------
$wsDeath = $wb.Worksheets.Item("Death")

$wsDeath.Range("A42").Value = 43938
$wsDeath.Range("A42").NumberFormat = "d-mmm-yy"

$wsDeath.Range("C42").Value = 15
$wsDeath.Range("C42").HorizontalAlignment = -4108
$wsDeath.Range("C42").VerticalAlignment = -4108

$wsDeath.Range("B42").Formula = "=SUM(B41+C42)"
$wsDeath.Range("B42").HorizontalAlignment = -4108
$wsDeath.Range("B42").VerticalAlignment = -4108

# Death's selection moves too, but it must not end up as the active tab.
$wsDeath.Range("B51:B52").Select()

# ---------------------------------------------------------------------
# Recoverd (sheet2) - new row 42; ends up the active tab/sheet, so it is
# handled last.
# ---------------------------------------------------------------------
$wsRecoverd = $wb.Worksheets.Item("Recoverd")
$wsRecoverd.Activate()

$wsRecoverd.Range("A42").Value = 43938
$wsRecoverd.Range("A42").NumberFormat = "d-mmm-yy"

$wsRecoverd.Range("C42").Value = 9
$wsRecoverd.Range("C42").HorizontalAlignment = -4108
$wsRecoverd.Range("C42").VerticalAlignment = -4108

# B36:B41 were plain (non-shared) formulas; rewriting the whole run lets
# the engine fold them into one shared-formula group through B42.
$wsRecoverd.Range("B36:B42").Formula = "=SUM(B35+C36)"
$wsRecoverd.Range("B42").HorizontalAlignment = -4108
$wsRecoverd.Range("B42").VerticalAlignment = -4108

$wsRecoverd.Range("D42").Select()
